$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accuracy values for B2:B118
$values = @(
    0.71875,0.515625,0.4375,0.3125,0.3125,0.375,0.375,0.359375,0.359375,0.34375,
    0.328125,0.390625,0.3125,0.265625,0.234375,0.171875,0.203125,0.1875,0.203125,0.234375,
    0.28125,0.265625,0.25,0.234375,0.234375,0.234375,0.25,0.25,0.25,0.25,
    0.25,0.25,0.25,0.25,0.25,0.25,0.25,0.25,0.25,0.265625,
    0.28125,0.28125,0.28125,0.296875,0.296875,0.296875,0.28125,0.28125,0.28125,0.28125,
    0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,
    0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,
    0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,
    0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,
    0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,
    0.28125,0.203125,0.1875,0.21875,0.28125,0.1875,0.265625,0.34375,0.125,0.234375,
    0.203125,0.234375,0.21875,0.234375,0.203125,0.234375,0.2131147540983606
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Update the repr() object address shown in column A for rows 102-118
$oldTag = "<__main__.DisplayOutputs object at 0x7f612057f4c0>"
$newTag = "<__main__.DisplayOutputs object at 0x7f09b90b06a0>"
for ($row = 102; $row -le 118; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value -eq $oldTag) {
        $cell.Value = $newTag
    }
}

# Update the active selection as recorded in the sheet view.
# The final user action was a "Select All" (Ctrl+A) performed while the
# cursor was resting on P11; reproduce it by touching P11 first and then
# extending the selection to the full grid, which is what Excel's
# Select-All does.
$ws.Range("P11").Select()
$ws.Cells.Select()
$ws.Range("A1:XFD1048576").Select()
